$d = $word.ActiveDocument

# 1. Insert the "Vi förväntar oss..." paragraph right after the
#    "Nedan presenteras fynd..." paragraph (near the top of the document).
$introRange = $d.Content
$introRange.Find.Execute("Nedan presenteras fynd av naturvårdsarter och fridlysta arter som gjorts i det avverkningsanmälda området, samt relevanta utdrag ur standarderna för FSC, Chain of Custody, Controlled Wood och PEFC.") | Out-Null
$introPara = $introRange.Paragraphs(1)
$introPara.Range.InsertParagraphAfter()
$newPara = $introPara.Next()
$newPara.Range.Text = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."

# 2. Remove the old trailing copy of that paragraph together with the two
#    empty paragraphs that preceded it at the end of the document body.
$count = $d.Paragraphs.Count
$pStart = $d.Paragraphs($count - 2)
$pEnd = $d.Paragraphs($count)
$trailRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$trailRange.Delete()

# 3. Update the date shown in the "first page" header from 2023-11-13 to
#    2023-11-14.
$sec = $d.Sections.First
$hdr = $sec.Headers(2)
$hdr.Range.Find.Execute("2023-11-13", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-11-14", 2)
